# Refatoração do código para permitir pesquisas com mais de 5 alternativas
# e permitir a escolha do delimitador da planilha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Change the delimiter used inside the "Options" column (D) from "," to ";"
#        for the existing survey rows (2-4), since commas inside option lists
#        conflicted with using comma as a separator.
$ws.Range("D2").Value = "Chocolate ao Leite; Chocolate Amargo"
$ws.Range("D3").Value = "Humanas; Exatas"
$ws.Range("D4").Value = "Bolo, Sorverte; Chocolate; Pé de moleque"

# --- 2. Add a new survey row (row 5) demonstrating support for more than
#        5 alternatives (16 numeric options).
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

$ws.Range("A5").Value = "Qual seu número favorito"
$ws.Range("B5").Value = "números !!!!!"
$ws.Range("C5").Value = "Professor de matemática"
$ws.Range("D5").Value = "1; 2; 3; 4; 5; 6; 7; 8; 9; 10; 11; 12; 13; 14; 15; 16"
$ws.Range("E5").Value = "Sim"
$ws.Range("F5").Value = "Sim"
$ws.Range("G5").Value = "Sim"
$ws.Range("H5").Value = "Sim"

# Match the row height used by the other data rows (wrapped, taller rows).
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight

# --- 3. Update the selection shown when the workbook is reopened.
$ws.Range("G10").Select()
